# updated legacy GSC export data
# Rolling 90-day window: drop the oldest date row, shift everything up by
# one day/row, and append the new day's values at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$firstDataRow = 2
$lastDataRow  = 91

function Add-OneDay {
    param([string]$dateStr)

    $parts = $dateStr.Split("-")
    $y = [int]$parts[0]
    $m = [int]$parts[1]
    $d = [int]$parts[2]

    $daysInMonth = @(31,28,31,30,31,30,31,31,30,31,30,31)
    $isLeap = (($y % 4 -eq 0) -and ($y % 100 -ne 0)) -or ($y % 400 -eq 0)
    if ($isLeap) { $daysInMonth[1] = 29 }

    $d = $d + 1
    if ($d -gt $daysInMonth[$m - 1]) {
        $d = 1
        $m = $m + 1
        if ($m -gt 12) {
            $m = 1
            $y = $y + 1
        }
    }

    $ms = "$m"
    if ($m -lt 10) { $ms = "0$m" }
    $ds = "$d"
    if ($d -lt 10) { $ds = "0$d" }

    return "$y-$ms-$ds"
}

# 1. Snapshot the current (pre-edit) Date text and Https-URLs numbers for
#    every row in the table before we start overwriting anything.
$oldDate = @{}
$oldHttps = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $oldDate[$r]  = $ws.Cells.Item($r, 1).Text
    $oldHttps[$r] = [double]$ws.Cells.Item($r, 3).Value2
}

# The brand-new day appended at the end of the window is one day after the
# last date that was already present.
$newLastDate = Add-OneDay $oldDate[$lastDataRow]

# 2. Shift column A (Date) and column C (HTTPS URLs) up by one row: row r
#    takes on what used to live in row r+1. The final row gets the new date
#    (its HTTPS URLs count stays whatever it already was - no new data
#    point is being reported for it yet).
for ($r = $firstDataRow; $r -lt $lastDataRow; $r++) {
    $ws.Cells.Item($r, 1).Value = "'" + $oldDate[$r + 1]
    $ws.Cells.Item($r, 3).Value = $oldHttps[$r + 1]
}
$ws.Cells.Item($lastDataRow, 1).Value = "'" + $newLastDate
$ws.Cells.Item($lastDataRow, 3).Value = $oldHttps[$lastDataRow]

# 3. The leading apostrophes above force text-entry (so Excel doesn't
#    reinterpret the date-shaped strings as serial date numbers), which
#    leaves a quote-prefixed number format behind. Re-flatten column A back
#    to the sheet's plain/general formatting so styles are untouched.
$blank = $ws.Range("ZZ1")
$blank.Copy()
$ws.Range($ws.Cells.Item($firstDataRow, 1), $ws.Cells.Item($lastDataRow, 1)).PasteSpecial(-4122)
$excel.CutCopyMode = $false
